$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-26 first (while current row numbering still matches original layout)
$ws.Rows("8:26").Delete()

# Overwrite rows 2-7 with the merged "(name, [texts...])" string values
$ws.Range("A2").Value = "('Celestine Reef', ['Plane — Luvion', 'Creatures without flying or islandwalk can’t attack.', 'Whenever you roll {CHAOS}, until a player planeswalks, you can’t lose the game and your opponents can’t win the game.'])"
$ws.Range("A3").Value = "('Horizon Boughs', ['Plane — Pyrulea', 'All permanents untap during each player’s untap step.', 'Whenever you roll {CHAOS}, you may search your library for up to three basic land cards, put them onto the battlefield tapped, then shuffle your library.'])"
$ws.Range("A4").Value = "('Mirrored Depths', ['Plane — Karsus', 'Whenever a player casts a spell, that player flips a coin. If the player loses the flip, counter that spell.', 'Whenever you roll {CHAOS}, target player reveals the top card of their library. If it’s a nonland card, you may cast it without paying its mana cost.'])"
$ws.Range("A5").Value = "('Stairs to Infinity', ['Plane — Xerex', 'Players have no maximum hand size.', 'Whenever you roll the planar die, draw a card.', 'Whenever you roll {CHAOS}, reveal the top card of your planar deck. You may put it on the bottom of your planar deck.'])"
$ws.Range("A6").Value = "('Tazeem', ['Plane — Zendikar', 'Creatures can’t block.', 'Whenever you roll {CHAOS}, draw a card for each land you control.'])"
$ws.Range("A7").Value = "('Tember City', ['Plane — Kinshala', 'Whenever a player taps a land for mana, Tember City deals 1 damage to that player.', 'Whenever you roll {CHAOS}, each other player sacrifices a nonland permanent.'])"
